$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# 1. tb_n_contact point estimate 7 -> 12
$ws.Cells.Item(2, 2).Value = 12

# 2. Insert three new rows before the old row 4 (susceptible_fully), shifting it (and "active") down
$ws.Range("A4:A6").EntireRow.Insert()

# Row 4: tb_prop_amplification
$ws.Cells.Item(4, 1).Value = "tb_prop_amplification"
$ws.Cells.Item(4, 2).Value = 0.8

# Row 5: start_mdr_introduce_time
$ws.Cells.Item(5, 1).Value = "start_mdr_introduce_time"
$ws.Cells.Item(5, 2).Value = 1880
$ws.Cells.Item(5, 5).Value = "Calendar year that MDR-TB first begins to emerge"

# Row 6: end_mdr_introduce_time
$ws.Cells.Item(6, 1).Value = "end_mdr_introduce_time"
$ws.Cells.Item(6, 2).Value = 1885
$ws.Cells.Item(6, 5).Value = "Calendar year that MDR-TB amplification reaches full parameter value"

# 3. susceptible_fully (now row 7) point estimate 3410000 -> 3200000
$ws.Cells.Item(7, 2).Value = 3200000

# 4. New row 9: age_breakpoints, with lower/upper-style values in B and C
$ws.Cells.Item(9, 1).Value = "age_breakpoints"
$ws.Cells.Item(9, 2).Value = 5
$ws.Cells.Item(9, 3).Value = 15

# New data validation for the two new numeric-year rows (decimal, -10000..10000)
$ws.Range("B5:D6").Validation.Add(2, 1, 1, -10000, 10000)

# Selection as in the authored workbook
$ws.Range("B7").Select()
